$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 17.01.2022 01:30"

# Makro row (row 5): D5 delta becomes a real number, E5 becomes a real date serial
$ws.Range("D5").Value = 0.6
$ws.Range("E5").Value = 44578.05217592593
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
